# Update cryptos list with latest prices / volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ col = value }
$updates = @{
    2  = @{ D = "34.176.79";  E = "  +0.52%  " }
    3  = @{ D = "1.788.27";   E = "  +0.38%  " }
    4  = @{ E = "  +0.15%  " }
    5  = @{ D = "226.66";     E = "  -0.29%  " }
    6  = @{ E = "  -0.36%  " }
    7  = @{ E = "  +0.13%  " }
    8  = @{ D = "32.07";      E = "  -1.29%  " }
    9  = @{ E = "  +2.91%  " }
    10 = @{ D = "0.0692";     E = "  -2.89%  " }
    11 = @{ D = "0.0945";     E = "  +1.00%  " }
    12 = @{ D = "2.047.13";   E = "  +0.45%  " }
    13 = @{ D = "11.36";      E = "  +1.86%  " }
    14 = @{ D = "1.789.64";   E = "  +0.93%  " }
    15 = @{ D = "34.132.25";  E = "  +0.39%  " }
    16 = @{ E = "  +0.51%  " }
    17 = @{ E = "  +2.20%  " }
    18 = @{ D = "68.02";      E = "  +0.39%  " }
    19 = @{ D = "245.85";     E = "  +0.59%  " }
    20 = @{ E = "  -0.33%  " }
    21 = @{ B = "Dai";        C = "https://coinranking.com/coin/MoTuySvg7+dai-dai";              D = "1.00";   E = "  +0.04%  " }
    22 = @{ B = "Avalanche";  C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax";        D = "10.86";  E = "  +1.49%  " }
    23 = @{ D = "4.11" }
    24 = @{ D = "2.05";       E = "  -1.18%  " }
    25 = @{ D = "161.58";     E = "  +1.06%  " }
    26 = @{ D = "7.15";       E = "  +1.46%  " }
    27 = @{ D = "16.32";      E = "  +0.26%  " }
    28 = @{ E = "  +1.70%  " }
    29 = @{ E = "  +0.30%  " }
    30 = @{ E = "  +0.32%  " }
    31 = @{ D = "0.0520";     E = "  +1.90%  " }
    32 = @{ E = "  +1.18%  " }
    33 = @{ D = "3.61";       E = "  +3.45%  " }
    34 = @{ E = "  +1.61%  " }
    35 = @{ D = "1.440.09";   E = "  +3.68%  " }
    36 = @{ D = "0.649";      E = "  -0.08%  " }
    37 = @{ E = "  +11.09%  " }
    38 = @{ E = "  +2.86%  " }
    39 = @{ E = "  -0.54%  " }
    40 = @{ D = "80.25";      E = "  +3.26%  " }
    41 = @{ E = "  +0.45%  " }
    42 = @{ D = "0.921";      E = "  +1.12%  " }
    43 = @{ D = "2.70";       E = "  +0.47%  " }
    44 = @{ D = "13.38";      E = "  +3.89%  " }
    45 = @{ B = "BabyDogeCoin"; C = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D = "0.0₆0139"; E = "  +0.87%  " }
    46 = @{ B = "Kaspa";      C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas";              D = "0.0511"; E = "  +2.75%  " }
    47 = @{ E = "  +4.81%  " }
    48 = @{ E = "  -0.45%  " }
    49 = @{ D = "107.74";     E = "  +0.17%  " }
    50 = @{ D = "1.948.71";   E = "  +0.68%  " }
    51 = @{ E = "  +0.15%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        if ($col -eq "D") {
            # The Price column holds values like "34.176.79" / "0.0520" /
            # "1.00" that must stay plain text (matching the original
            # inline strings) instead of being auto-coerced into numbers
            # by Excel's smart-entry logic (which would e.g. turn
            # "0.0520" into 0.052 or "1.00" into 1). Force Text format on
            # just this cell before assigning, same as it already is in
            # the source file.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $cols[$col]
    }
}
